{"js": "const replacements = [\n  { oldText: \"A imagem mostra a interface do aplicativo Google Play, especificamente a se\u00e7\u00e3o \\\"Para voc\u00ea\\\", onde s\u00e3o exibidos jogos e sugest\u00f5es baseadas na atividade recente do usu\u00e1rio. Na parte superior, destacam-se jogos como \\\"Blood Strike\\\", \\\"Roblox\\\" e \\\"Tile Club\\\". Abaixo, a se\u00e7\u00e3o de \\\"Sugest\u00f5es para voc\u00ea\\\" apresenta jogos patrocinados, incluindo \\\"Paci\u00eancia - Jogo de Solit\u00e1rio\\\", \\\"Coin Master\\\" e \\\"Bubble Pop! Cannon Shooter\\\". O canto inferior direito cont\u00e9m um \u00edcone de pesquisa destacado em vermelho.\", newText: \"A imagem apresenta a tela inicial da Google Play Store, mostrando uma se\u00e7\u00e3o personalizada chamada \\\"Para voc\u00ea\\\". Ela lista sugest\u00f5es de jogos baseadas na atividade recente do usu\u00e1rio, incluindo t\u00edtulos como \\\"Blood Strike\\\", \\\"Roblox\\\" e \\\"Tile Club\\\". Abaixo, h\u00e1 uma se\u00e7\u00e3o de jogos patrocinados, com sugest\u00f5es como \\\"Paci\u00eancia - Jogo de Solit\u00e1rio\\\", \\\"Coin Master\\\" e \\\"Bubble Pop! Cannon Shooter\\\". A parte inferior da tela exibe \u00edcones para navegar entre Jogos, Apps e Livros, com uma barra de pesquisa destacada em vermelho.\" },\n  { oldText: \"A imagem mostra a tela de busca em um aplicativo de loja (provavelmente Google Play) com o termo \\\"smart sales force\\\" inserido. Na parte superior, h\u00e1 um t\u00edtulo destacado. A listagem apresenta v\u00e1rios aplicativos relacionados, incluindo \\\"Smart For\u00e7a de Vendas\\\" da Arpa Sistemas, que tem uma avalia\u00e7\u00e3o de 4,3 estrelas e requer 14 MB de armazenamento. Os aplicativos s\u00e3o classificados por popularidade e tipo, com informa\u00e7\u00f5es adicionais como n\u00famero de avalia\u00e7\u00f5es e categorias.\", newText: \"A imagem mostra uma tela de pesquisa em um aplicativo de loja digital, onde o usu\u00e1rio procura por \\\"smart sales force\\\". Na parte superior, aparecem os resultados patrocinados relacionados a essa busca. O destaque \u00e9 para o aplicativo \\\"Smart For\u00e7a de Vendas\\\" da Arpa Sistemas, que possui uma avalia\u00e7\u00e3o de 4,3 estrelas, ocupa 14 MB de espa\u00e7o e tem mais de mil downloads. Outros aplicativos relacionados, como Salesforce e App Sales Force, tamb\u00e9m est\u00e3o listados abaixo.\" },\n  { oldText: \"A imagem apresenta a p\u00e1gina de download do aplicativo \\\"Smart For\u00e7a de Vendas\\\" na Google Play Store. O aplicativo \u00e9 desenvolvido pela Arpa Sistemas e possui uma classifica\u00e7\u00e3o de 4,2 estrelas com 12 avalia\u00e7\u00f5es. O tamanho do arquivo \u00e9 de 14 MB. Na parte superior, h\u00e1 o \u00edcone do aplicativo e, abaixo, algumas capturas de tela que mostram a interface do usu\u00e1rio. Tamb\u00e9m est\u00e3o dispon\u00edveis informa\u00e7\u00f5es sobre o aplicativo e op\u00e7\u00f5es relacionadas a \\\"Neg\u00f3cios\\\" e \\\"Seguran\u00e7a dos dados\\\". Um bot\u00e3o destacado no centro permite a instala\u00e7\u00e3o do aplicativo.\", newText: \"A imagem apresenta a p\u00e1gina de download do aplicativo \\\"Smart For\u00e7a de Vendas\\\" na Google Play Store. O aplicativo, desenvolvido pela Arpa Sistemas, possui uma classifica\u00e7\u00e3o de 4,2 estrelas, com 12 avalia\u00e7\u00f5es e um tamanho de 14 MB. A interface exibe capturas de tela do aplicativo em uso, mostrando diferentes funcionalidades e op\u00e7\u00f5es. H\u00e1 uma se\u00e7\u00e3o informativa sobre o aplicativo, destacando que \\\"Smart Vendas\\\" agora \u00e9 chamado de \\\"Smart For\u00e7a de Vendas\\\", e um bot\u00e3o para instal\u00e1-lo. Na parte inferior da tela, h\u00e1 \u00edcones que permitem acessar jogos, outros aplicativos, e livros, al\u00e9m de uma op\u00e7\u00e3o de pesquisa.\" },\n  { oldText: \"A imagem mostra uma tela de instala\u00e7\u00e3o do aplicativo \\\"Smart For\u00e7a de Vendas\\\" em um dispositivo Android. No topo, h\u00e1 a indica\u00e7\u00e3o de que o aplicativo est\u00e1 sendo instalado e que \u00e9 verificado pelo Play Protect. Abaixo, est\u00e3o sugest\u00f5es de aplicativos patrocinados, incluindo \u00edcones de apps como \\\"Nomad\\\", \\\"Livelo\\\", e \\\"Estoque, Vendas, Pdv, Finan\u00e7as\\\", com suas respectivas classifica\u00e7\u00f5es. Na parte inferior, h\u00e1 op\u00e7\u00f5es para acesso a mais aplicativos e a barra de navega\u00e7\u00e3o com \u00edcones de jogos, aplicativos, e livros.\", newText: \"A imagem mostra a tela de instala\u00e7\u00e3o do aplicativo \\\"Smart For\u00e7a de Vendas\\\" em um dispositivo m\u00f3vel. Na parte superior, h\u00e1 uma barra de status com o hor\u00e1rio e a qualidade do sinal. Abaixo, s\u00e3o apresentadas sugest\u00f5es de aplicativos patrocinados, como \\\"Nomad: Conta em D\u00f3lar e Cart\u00e3o\\\", \\\"Livelo: juntar e trocar pontos\\\" e \\\"Estoque, Vendas, PDV, Finan\u00e7as\\\". Na parte inferior, h\u00e1 uma se\u00e7\u00e3o com mais aplicativos para testar, incluindo \\\"PictureThis\\\", \\\"Arquivos do Google\\\" e \\\"CamScanner\\\". A interface \u00e9 simples, com \u00edcones coloridos representando cada aplicativo.\" },\n];\n\nconst body = context.document.body;\nlet appliedCount = 0;\n\nfor (const { oldText, newText } of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      `Expected exactly 1 match, found ${results.items.length} for: ` + oldText.substring(0, 60)\n    );\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  appliedCount++;\n  await context.sync();\n}\n\nreturn `replaced ${appliedCount} paragraph(s)`;\n", "ps1": "$d = $word.ActiveDocument\n$replacements = @(\n    @{ Old = 'A imagem mostra a interface do aplicativo Google Play, especificamente a se\u00e7\u00e3o \"Para voc\u00ea\", onde s\u00e3o exibidos jogos e sugest\u00f5es baseadas na atividade recente do usu\u00e1rio. Na parte superior, destacam-se jogos como \"Blood Strike\", \"Roblox\" e \"Tile Club\". Abaixo, a se\u00e7\u00e3o de \"Sugest\u00f5es para voc\u00ea\" apresenta jogos patrocinados, incluindo \"Paci\u00eancia - Jogo de Solit\u00e1rio\", \"Coin Master\" e \"Bubble Pop! Cannon Shooter\". O canto inferior direito cont\u00e9m um \u00edcone de pesquisa destacado em vermelho.'; New = 'A imagem apresenta a tela inicial da Google Play Store, mostrando uma se\u00e7\u00e3o personalizada chamada \"Para voc\u00ea\". Ela lista sugest\u00f5es de jogos baseadas na atividade recente do usu\u00e1rio, incluindo t\u00edtulos como \"Blood Strike\", \"Roblox\" e \"Tile Club\". Abaixo, h\u00e1 uma se\u00e7\u00e3o de jogos patrocinados, com sugest\u00f5es como \"Paci\u00eancia - Jogo de Solit\u00e1rio\", \"Coin Master\" e \"Bubble Pop! Cannon Shooter\". A parte inferior da tela exibe \u00edcones para navegar entre Jogos, Apps e Livros, com uma barra de pesquisa destacada em vermelho.' },\n    @{ Old = 'A imagem mostra a tela de busca em um aplicativo de loja (provavelmente Google Play) com o termo \"smart sales force\" inserido. Na parte superior, h\u00e1 um t\u00edtulo destacado. A listagem apresenta v\u00e1rios aplicativos relacionados, incluindo \"Smart For\u00e7a de Vendas\" da Arpa Sistemas, que tem uma avalia\u00e7\u00e3o de 4,3 estrelas e requer 14 MB de armazenamento. Os aplicativos s\u00e3o classificados por popularidade e tipo, com informa\u00e7\u00f5es adicionais como n\u00famero de avalia\u00e7\u00f5es e categorias.'; New = 'A imagem mostra uma tela de pesquisa em um aplicativo de loja digital, onde o usu\u00e1rio procura por \"smart sales force\". Na parte superior, aparecem os resultados patrocinados relacionados a essa busca. O destaque \u00e9 para o aplicativo \"Smart For\u00e7a de Vendas\" da Arpa Sistemas, que possui uma avalia\u00e7\u00e3o de 4,3 estrelas, ocupa 14 MB de espa\u00e7o e tem mais de mil downloads. Outros aplicativos relacionados, como Salesforce e App Sales Force, tamb\u00e9m est\u00e3o listados abaixo.' },\n    @{ Old = 'A imagem apresenta a p\u00e1gina de download do aplicativo \"Smart For\u00e7a de Vendas\" na Google Play Store. O aplicativo \u00e9 desenvolvido pela Arpa Sistemas e possui uma classifica\u00e7\u00e3o de 4,2 estrelas com 12 avalia\u00e7\u00f5es. O tamanho do arquivo \u00e9 de 14 MB. Na parte superior, h\u00e1 o \u00edcone do aplicativo e, abaixo, algumas capturas de tela que mostram a interface do usu\u00e1rio. Tamb\u00e9m est\u00e3o dispon\u00edveis informa\u00e7\u00f5es sobre o aplicativo e op\u00e7\u00f5es relacionadas a \"Neg\u00f3cios\" e \"Seguran\u00e7a dos dados\". Um bot\u00e3o destacado no centro permite a instala\u00e7\u00e3o do aplicativo.'; New = 'A imagem apresenta a p\u00e1gina de download do aplicativo \"Smart For\u00e7a de Vendas\" na Google Play Store. O aplicativo, desenvolvido pela Arpa Sistemas, possui uma classifica\u00e7\u00e3o de 4,2 estrelas, com 12 avalia\u00e7\u00f5es e um tamanho de 14 MB. A interface exibe capturas de tela do aplicativo em uso, mostrando diferentes funcionalidades e op\u00e7\u00f5es. H\u00e1 uma se\u00e7\u00e3o informativa sobre o aplicativo, destacando que \"Smart Vendas\" agora \u00e9 chamado de \"Smart For\u00e7a de Vendas\", e um bot\u00e3o para instal\u00e1-lo. Na parte inferior da tela, h\u00e1 \u00edcones que permitem acessar jogos, outros aplicativos, e livros, al\u00e9m de uma op\u00e7\u00e3o de pesquisa.' },\n    @{ Old = 'A imagem mostra uma tela de instala\u00e7\u00e3o do aplicativo \"Smart For\u00e7a de Vendas\" em um dispositivo Android. No topo, h\u00e1 a indica\u00e7\u00e3o de que o aplicativo est\u00e1 sendo instalado e que \u00e9 verificado pelo Play Protect. Abaixo, est\u00e3o sugest\u00f5es de aplicativos patrocinados, incluindo \u00edcones de apps como \"Nomad\", \"Livelo\", e \"Estoque, Vendas, Pdv, Finan\u00e7as\", com suas respectivas classifica\u00e7\u00f5es. Na parte inferior, h\u00e1 op\u00e7\u00f5es para acesso a mais aplicativos e a barra de navega\u00e7\u00e3o com \u00edcones de jogos, aplicativos, e livros.'; New = 'A imagem mostra a tela de instala\u00e7\u00e3o do aplicativo \"Smart For\u00e7a de Vendas\" em um dispositivo m\u00f3vel. Na parte superior, h\u00e1 uma barra de status com o hor\u00e1rio e a qualidade do sinal. Abaixo, s\u00e3o apresentadas sugest\u00f5es de aplicativos patrocinados, como \"Nomad: Conta em D\u00f3lar e Cart\u00e3o\", \"Livelo: juntar e trocar pontos\" e \"Estoque, Vendas, PDV, Finan\u00e7as\". Na parte inferior, h\u00e1 uma se\u00e7\u00e3o com mais aplicativos para testar, incluindo \"PictureThis\", \"Arquivos do Google\" e \"CamScanner\". A interface \u00e9 simples, com \u00edcones coloridos representando cada aplicativo.' },\n)\n\n$totalReplaced = 0\nforeach ($pair in $replacements) {\n    $done = $false\n    foreach ($p in $d.Paragraphs) {\n        $t = $p.Range.Text\n        if ($t.TrimEnd([char]13) -eq $pair.Old) {\n            $p.Range.Text = $pair.New\n            $totalReplaced++\n            $done = $true\n            break\n        }\n    }\n    if (-not $done) {\n        throw \"No paragraph matched target text starting with: \" + $pair.Old.Substring(0, 40)\n    }\n}\n\nWrite-Output \"replaced $totalReplaced paragraph(s)\"\n"}
